# Auto-generated script to apply scheduled runner updates to Excalibur_Profits workbook
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ derived columns
# (H, I, J, K, L, M, N) across several rows on multiple sheets, reflecting refreshed market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 27741.398
$ws.Range("I132").Value = 30587.367
$ws.Range("K132").Value = 91762.101
$ws.Range("M132").Value = -89232.101

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2953.9697
$ws.Range("I138").Value = 2357.4375
$ws.Range("J138").Value = 3515.4119
$ws.Range("K138").Value = 7072.3125
$ws.Range("L138").Value = 10546.2357
$ws.Range("M138").Value = -1932.3125
$ws.Range("N138").Value = -20826.2357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2980.5
$ws.Range("I2").Value = 3077.7856
$ws.Range("K2").Value = 3077.7856
$ws.Range("M2").Value = -2964.7856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3479.87
$ws.Range("I32").Value = 1586.1097
$ws.Range("K32").Value = 1586.1097
$ws.Range("M32").Value = -1299.1097

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6585.0435
$ws.Range("I45").Value = 5098.875
$ws.Range("K45").Value = 5098.875
$ws.Range("M45").Value = -4721.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50636

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1607.1333
$ws.Range("I97").Value = 1257.6428
$ws.Range("J97").Value = 6500
$ws.Range("K97").Value = 1257.6428
$ws.Range("L97").Value = 6500
$ws.Range("M97").Value = -761.6428000000001
$ws.Range("N97").Value = -7492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1829.7576
$ws.Range("I102").Value = 1367.4839
$ws.Range("K102").Value = 1367.4839
$ws.Range("M102").Value = 254.5161000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 40387
$ws.Range("J112").Value = 40387
$ws.Range("L112").Value = 40387
$ws.Range("N112").Value = -43341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2980.5
$ws.Range("I116").Value = 3077.7856
$ws.Range("K116").Value = 3077.7856
$ws.Range("M116").Value = -783.7856000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 200279.83
$ws.Range("J119").Value = 200279.83
$ws.Range("L119").Value = 200279.83
$ws.Range("N119").Value = -209955.83

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3035.3794
$ws.Range("I122").Value = 1777.25
$ws.Range("J122").Value = 3923.4707
$ws.Range("K122").Value = 5331.75
$ws.Range("L122").Value = 11770.4121
$ws.Range("M122").Value = -2881.75
$ws.Range("N122").Value = -16670.4121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2980.5
$ws.Range("I3").Value = 3077.7856
$ws.Range("K3").Value = 3077.7856
$ws.Range("M3").Value = -2963.7856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 74107.664
$ws.Range("J50").Value = 74107.664
$ws.Range("L50").Value = 74107.664
$ws.Range("N50").Value = -75255.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2320.3076
$ws.Range("I64").Value = 355.6
$ws.Range("J64").Value = 3548.25
$ws.Range("K64").Value = 355.6
$ws.Range("L64").Value = 3548.25
$ws.Range("M64").Value = -130.6
$ws.Range("N64").Value = -3998.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 2320.3076
$ws.Range("I67").Value = 355.6
$ws.Range("J67").Value = 3548.25
$ws.Range("K67").Value = 355.6
$ws.Range("L67").Value = 3548.25
$ws.Range("M67").Value = 424.4
$ws.Range("N67").Value = -5108.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2477.5715
$ws.Range("I99").Value = 1668.8
$ws.Range("K99").Value = 1668.8
$ws.Range("M99").Value = -170.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2352.4512
$ws.Range("I134").Value = 1566.4084
$ws.Range("J134").Value = 7426
$ws.Range("K134").Value = 4699.2252
$ws.Range("L134").Value = 22278
$ws.Range("M134").Value = -2164.2252
$ws.Range("N134").Value = -27348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7226.4053
$ws.Range("I31").Value = 2512.4546
$ws.Range("J31").Value = 9220.77
$ws.Range("K31").Value = 2512.4546
$ws.Range("L31").Value = 9220.77
$ws.Range("M31").Value = -2217.4546
$ws.Range("N31").Value = -9810.77

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7226.4053
$ws.Range("I34").Value = 2512.4546
$ws.Range("J34").Value = 9220.77
$ws.Range("K34").Value = 2512.4546
$ws.Range("L34").Value = 9220.77
$ws.Range("M34").Value = -2310.4546
$ws.Range("N34").Value = -9624.77

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 25000
$ws.Range("J38").Value = 25000
$ws.Range("L38").Value = 25000
$ws.Range("N38").Value = -25754

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 25000
$ws.Range("J46").Value = 25000
$ws.Range("L46").Value = 25000
$ws.Range("N46").Value = -25422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5961.222
$ws.Range("I99").Value = 5392.9287
$ws.Range("K99").Value = 5392.9287
$ws.Range("M99").Value = -3894.9287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 59666.332
$ws.Range("J100").Value = 59666.332
$ws.Range("L100").Value = 59666.332
$ws.Range("N100").Value = -61830.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5961.222
$ws.Range("I126").Value = 5392.9287
$ws.Range("K126").Value = 16178.7861
$ws.Range("M126").Value = -13708.7861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 277818.12
$ws.Range("J2").Value = 41
$ws.Range("L2").Value = 246
$ws.Range("N2").Value = -472

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1750233.1
$ws.Range("I7").Value = 293.5
$ws.Range("K7").Value = 880.5
$ws.Range("M7").Value = -768.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 852.6923
$ws.Range("I14").Value = 852.6923
$ws.Range("K14").Value = 2558.0769
$ws.Range("M14").Value = -2385.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3350.6
$ws.Range("I70").Value = 2189.5
$ws.Range("K70").Value = 6568.5
$ws.Range("M70").Value = -6253.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3350.6
$ws.Range("I73").Value = 2189.5
$ws.Range("K73").Value = 6568.5
$ws.Range("M73").Value = -5476.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 862196.6
$ws.Range("I132").Value = 862196.6
$ws.Range("K132").Value = 2586589.8
$ws.Range("M132").Value = -2584059.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5952.75
$ws.Range("I7").Value = 5759.222
$ws.Range("K7").Value = 5759.222
$ws.Range("M7").Value = -5647.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3343.7058
$ws.Range("I40").Value = 3087
$ws.Range("J40").Value = 3959.8
$ws.Range("K40").Value = 3087
$ws.Range("L40").Value = 3959.8
$ws.Range("M40").Value = -2951
$ws.Range("N40").Value = -4231.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3558.3333
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 3633.9285
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 3633.9285
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -4009.9285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2583.6667
$ws.Range("J61").Value = 2124.5
$ws.Range("L61").Value = 2124.5
$ws.Range("N61").Value = -2528.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5561.5
$ws.Range("I100").Value = 2247.3635
$ws.Range("J100").Value = 14675.375
$ws.Range("K100").Value = 2247.3635
$ws.Range("L100").Value = 14675.375
$ws.Range("M100").Value = -1706.3635
$ws.Range("N100").Value = -15757.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2583.6667
$ws.Range("J113").Value = 2124.5
$ws.Range("L113").Value = 2124.5
$ws.Range("N113").Value = -6464.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4823.278
$ws.Range("I122").Value = 4099.6665
$ws.Range("K122").Value = 12298.9995
$ws.Range("M122").Value = -9848.999500000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5952.75
$ws.Range("I126").Value = 5759.222
$ws.Range("K126").Value = 17277.666
$ws.Range("M126").Value = -14807.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 167593660
$ws.Range("I4").Value = 2755500
$ws.Range("K4").Value = 2755500
$ws.Range("M4").Value = -2755387

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 120698
$ws.Range("J119").Value = 120698
$ws.Range("L119").Value = 120698
$ws.Range("N119").Value = -130374

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2696.5454
$ws.Range("I122").Value = 1958.5625
$ws.Range("K122").Value = 5875.6875
$ws.Range("M122").Value = -3425.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5465.8887
$ws.Range("I126").Value = 6239.4
$ws.Range("J126").Value = 4499
$ws.Range("K126").Value = 18718.2
$ws.Range("L126").Value = 13497
$ws.Range("M126").Value = -16248.2
$ws.Range("N126").Value = -18437
